# Updates betting-odds data on the "Poland I Liga" sheet.
# The underlying source rows got re-ordered/re-fetched upstream, which in the
# canonical OOXML shows up as whole rows (columns B:AC) being rotated among a
# handful of row groups, plus a few odds corrections applied in place on
# another group of rows. This script reproduces both kinds of edits using the
# Excel COM object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row group 1 (match ids 5140743 / 5139053 / 5139054): 3-way rotation.
#   row 51 <- old row 52
#   row 52 <- old row 53
#   row 53 <- old row 51
# ---------------------------------------------------------------------------
$row51 = $ws.Range("B51:AC51").Value2
$row52 = $ws.Range("B52:AC52").Value2
$row53 = $ws.Range("B53:AC53").Value2

$ws.Range("B51:AC51").Value2 = $row52
$ws.Range("B52:AC52").Value2 = $row53
$ws.Range("B53:AC53").Value2 = $row51

# ---------------------------------------------------------------------------
# Row group 2 (match ids around 5448048 / 5451608 / 5451609 / 5452381):
# 4-way rotation; row 138 is untouched.
#   row 136 <- old row 140
#   row 137 <- old row 139
#   row 139 <- old row 136
#   row 140 <- old row 137
# ---------------------------------------------------------------------------
$row136 = $ws.Range("B136:AC136").Value2
$row137 = $ws.Range("B137:AC137").Value2
$row139 = $ws.Range("B139:AC139").Value2
$row140 = $ws.Range("B140:AC140").Value2

$ws.Range("B136:AC136").Value2 = $row140
$ws.Range("B137:AC137").Value2 = $row139
$ws.Range("B139:AC139").Value2 = $row136
$ws.Range("B140:AC140").Value2 = $row137

# ---------------------------------------------------------------------------
# Row group 3 (match ids 5448049 / 5448050): simple swap.
# ---------------------------------------------------------------------------
$row143 = $ws.Range("B143:AC143").Value2
$row144 = $ws.Range("B144:AC144").Value2

$ws.Range("B143:AC143").Value2 = $row144
$ws.Range("B144:AC144").Value2 = $row143

# ---------------------------------------------------------------------------
# Row group 4 (match ids 6805719 / 6803727): simple swap.
# ---------------------------------------------------------------------------
$row209 = $ws.Range("B209:AC209").Value2
$row210 = $ws.Range("B210:AC210").Value2

$ws.Range("B209:AC209").Value2 = $row210
$ws.Range("B210:AC210").Value2 = $row209

# ---------------------------------------------------------------------------
# Rows 318-321: a handful of odds columns corrected in place (same match ids,
# no row re-ordering here).
# ---------------------------------------------------------------------------
$ws.Range("R318").Value2 = 1.9
$ws.Range("S318").Value2 = 1.95
$ws.Range("U318").Value2 = 1.875
$ws.Range("V318").Value2 = 1.975

$ws.Range("N319").Value2 = 4
$ws.Range("P319").Value2 = 1.909
$ws.Range("R319").Value2 = 1.95
$ws.Range("S319").Value2 = 1.9

$ws.Range("N320").Value2 = 2.2
$ws.Range("O320").Value2 = 3.25
$ws.Range("P320").Value2 = 3.3
$ws.Range("R320").Value2 = 1.975
$ws.Range("S320").Value2 = 1.875
$ws.Range("T320").Value2 = 2.25
$ws.Range("U320").Value2 = 1.85
$ws.Range("V320").Value2 = 2

$ws.Range("N321").Value2 = 4.333
$ws.Range("P321").Value2 = 1.75
$ws.Range("Q321").Value2 = 0.75
$ws.Range("R321").Value2 = 1.875
$ws.Range("S321").Value2 = 1.975
